$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: set a Price cell as text (avoids Excel auto-numeric coercion)
# by forcing NumberFormat to text, assigning the value, then clearing the
# transient format so the cell keeps its original (unstyled) appearance.

# Rows 2-41: update Price (D) and/or Volume(1h) (E) columns
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.119.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.462.42"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.640"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.27%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.69"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000270"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.035.38"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.473.67"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.296.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.986"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.80"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.09"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.53"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "613.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -10.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.36"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.35%  "
$ws.Range("E38").Value = "  +8.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0785"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.365.34"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.379"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.88%  "

# Rows 42-43: Stacks and FirstDigitalUSD swap positions with updated data
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.29%  "

# Rows 44-51: update Price (D) and/or Volume(1h) (E) columns
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -10.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0412"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.69"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.75%  "
